# Inserts a new "enable_year" parameter row (row 7) into the dataset,
# shifting all subsequent data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the current row 7 (the first data row after the
# "input" row); this shifts all rows from 7 downward by one, and Excel
# automatically extends the AutoFilter / used range / dimension.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row with the "enable_year" configuration entry.
$ws.Range("A7").Value = "CHE"
$ws.Range("B7").Value = "conv_elec_onshorewind"
$ws.Range("C7").Value = "enable_year"
$ws.Range("D7").Value = "configuration"
$ws.Range("G7").Value = 1990

# Re-apply the AutoFilter so its range grows to include the newly
# inserted row.
$ws.AutoFilterMode = $false
$ws.Range("A5:L853").AutoFilter()

# Excel keeps the hidden _xlnm._FilterDatabase name in sync with the
# AutoFilter range, but the COM layer doesn't always refresh it, so set
# it explicitly to match.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$5:`$L`$853"

# Restore the selection similar to the authored edit.
$ws.Range("D11").Select()

$wb.Save()
